# Update the "quotes" worksheet:
#   - replace the 9-row quote list with a new 7-row quote list (new quotes/authors)
#   - resize a couple of rows to fit the new text
#   - tighten the two data columns
#   - change body/header vertical alignment from "center" to "top"
#   - move the active selection to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the last three quote rows (rows 9-11) - the new list only
#    needs 7 quote rows (rows 2-8) under the header row.
# ------------------------------------------------------------------
$ws.Rows("9:11").Delete()

# ------------------------------------------------------------------
# 2. Write the new quotes/authors into A2:B8
# ------------------------------------------------------------------
$ws.Range("A2").Value = "“The more that you read, the more things you will know. The more that you learn, the more places you'll go.”"
$ws.Range("B2").Value = "Dr. Seuss"

$ws.Range("A3").Value = "“We read to know we're not alone.”"
$ws.Range("B3").Value = "William Nicholson"

$ws.Range("A4").Value = "“A reader lives a thousand lives before he dies, said Jojen. The man who never reads lives only one.”"
$ws.Range("B4").Value = "George R.R. Martin"

$ws.Range("A5").Value = "“You can never get a cup of tea large enough or a book long enough to suit me.”"
$ws.Range("B5").Value = "C.S. Lewis"

$ws.Range("A6").Value = "“What really knocks me out is a book that, when you're all done reading it, you wish the author that wrote it was a terrific friend of yours and you could call him up on the phone whenever you felt like it. That doesn't happen much, though.”"
$ws.Range("B6").Value = "J.D. Salinger"

$ws.Range("A7").Value = "“′Classic′ - a book which people praise and don't read.”"
$ws.Range("B7").Value = "Mark Twain"

$ws.Range("A8").Value = "“I declare after all there is no enjoyment like reading! How much sooner one tires of any thing than of a book! -- When I have a house of my own, I shall be miserable if I have not an excellent library.”"
$ws.Range("B8").Value = "Jane Austen"

# ------------------------------------------------------------------
# 3. Row heights for the new text lengths
# ------------------------------------------------------------------
$ws.Rows(2).RowHeight = 43.5
$ws.Rows(3).RowHeight = 43.5
$ws.Rows(4).RowHeight = 43.5
$ws.Rows(5).RowHeight = 43.5
$ws.Rows(6).RowHeight = 101.5
$ws.Rows(7).RowHeight = 29
$ws.Rows(8).RowHeight = 87

# ------------------------------------------------------------------
# 4. Narrower columns to match the new layout
# ------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 33.5
$ws.Columns("B").ColumnWidth = 13.16666667

# ------------------------------------------------------------------
# 5. Vertical alignment: center -> top for the header and the body
# ------------------------------------------------------------------
$ws.Range("A1:B8").VerticalAlignment = -4160

# ------------------------------------------------------------------
# 6. Move the selection to A2
# ------------------------------------------------------------------
$ws.Range("A2").Select()
